# Generate Report for Archive
#
# Two previously "In progress" files (3bf7a7e2-... and 99666f0c-...) have
# moved from "Ready for handoff" status into "In Translation" status, on
# every sheet of the report. The "b69fb0f1-..." row stays "Ready for
# handoff" on every sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("C4").Value = $newStatus

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("C4").Value = $newStatus
